$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D (Price) cells to Text format so that
# numeric-looking strings (e.g. "1.001") are not auto-converted to
# numbers by Excel when the new values are written.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "29.320.27"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "1.876.94"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "0.7104"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("D6").Value = "242.10"
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").Value = "0.07992"
$ws.Range("E8").Value = "  +2.52%  "
$ws.Range("D9").Value = "0.3152"
$ws.Range("E9").Value = "  +1.29%  "
$ws.Range("E10").Value = "  -0.67%  "
$ws.Range("D11").Value = "0.08275"
$ws.Range("E11").Value = "  -1.61%  "
$ws.Range("D12").Value = "1.889.15"
$ws.Range("E12").Value = "  +1.24%  "
$ws.Range("D13").Value = "5.243"
$ws.Range("E13").Value = "  +0.12%  "
$ws.Range("D14").Value = "94.40"
$ws.Range("E14").Value = "  +3.60%  "
$ws.Range("D15").Value = "0.7103"
$ws.Range("E15").Value = "  -0.24%  "
$ws.Range("D16").Value = "6.358"
$ws.Range("E16").Value = "  +4.36%  "
$ws.Range("D17").Value = "0.000008508"
$ws.Range("E17").Value = "  +3.40%  "
$ws.Range("D18").Value = "29.352.80"
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("D19").Value = "245.22"
$ws.Range("E19").Value = "  +2.16%  "
$ws.Range("D20").Value = "2.149.28"
$ws.Range("E20").Value = "  +1.70%  "
$ws.Range("D21").Value = "13.26"
$ws.Range("E21").Value = "  +0.39%  "
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("D23").Value = "7.774"
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("D24").Value = "1.004"
$ws.Range("E24").Value = "  +0.29%  "
$ws.Range("D25").Value = "0.1553"
$ws.Range("E25").Value = "  -2.90%  "
$ws.Range("D26").Value = "9.053"
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("D27").Value = "162.58"
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("D28").Value = "18.53"
$ws.Range("E28").Value = "  +0.23%  "
$ws.Range("D29").Value = "1.505"
$ws.Range("E29").Value = "  -0.29%  "
$ws.Range("D30").Value = "4.413"
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("D31").Value = "4.319"
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("D32").Value = "1.185"
$ws.Range("E32").Value = "  -8.06%  "
$ws.Range("D33").Value = "0.05369"
$ws.Range("E33").Value = "  +1.40%  "
$ws.Range("D35").Value = "0.7636"
$ws.Range("E35").Value = "  +2.79%  "
$ws.Range("D36").Value = "1.182"
$ws.Range("E36").Value = "  +0.43%  "
$ws.Range("D37").Value = "2.687"
$ws.Range("E37").Value = "  -0.48%  "
$ws.Range("D38").Value = "0.01881"
$ws.Range("E38").Value = "  +0.46%  "
$ws.Range("D39").Value = "1.257.27"
$ws.Range("E39").Value = "  +3.05%  "
$ws.Range("D40").Value = "2.751"
$ws.Range("E40").Value = "  +0.83%  "
$ws.Range("D41").Value = "6.513"
$ws.Range("E41").Value = "  -0.60%  "
$ws.Range("D42").Value = "112.75"
$ws.Range("E42").Value = "  +1.71%  "
$ws.Range("D43").Value = "0.9125"
$ws.Range("E43").Value = "  +2.91%  "
$ws.Range("D44").Value = "74.11"
$ws.Range("E44").Value = "  +1.98%  "
$ws.Range("D45").Value = "0.00000000132"
$ws.Range("E45").Value = "  +8.52%  "
$ws.Range("D46").Value = "1.001"
$ws.Range("E46").Value = "  +0.12%  "
$ws.Range("D47").Value = "2.035.85"
$ws.Range("E47").Value = "  +0.83%  "
$ws.Range("D48").Value = "0.5224"
$ws.Range("E48").Value = "  +0.74%  "
$ws.Range("D49").Value = "1.797"
$ws.Range("E49").Value = "  -0.30%  "
$ws.Range("D50").Value = "9.443"
$ws.Range("E50").Value = "  +0.46%  "
$ws.Range("D51").Value = "0.4367"
$ws.Range("E51").Value = "  +1.13%  "

# Restore the original (default) cell style on column D now that the
# text values are safely stored, so formatting matches the source file.
$priceRange.Style = "Normal"
